$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 95 with the missing time-log entry
$ws.Range("A95").Value = 41932
$ws.Range("D95").Value = 5
$ws.Range("B95").Value = 0.93055555555555547
$ws.Range("C95").Value = 0.99930555555555556
$ws.Range("F95").Value = "Coding"

# Update the visible window / selection to match where the user ended up
$ws.Application.ActiveWindow.ScrollRow = 76
$ws.Range("B96").Select()
